# Apply the Tonberry_Profits scheduled-runner update across all 8 sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 499.5
$ws.Range("I28").Value = 333.33334
$ws.Range("K28").Value = 333.33334
$ws.Range("M28").Value = 151.66666
$ws.Range("H32").Value = 1250
$ws.Range("J32").Value = 1500
$ws.Range("L32").Value = 1500
$ws.Range("N32").Value = -2152
$ws.Range("H33").Value = 270.1905
$ws.Range("I33").Value = 410.54544
$ws.Range("J33").Value = 115.8
$ws.Range("K33").Value = 410.54544
$ws.Range("L33").Value = 115.8
$ws.Range("M33").Value = -181.54544
$ws.Range("N33").Value = -573.8
$ws.Range("H62").Value = 2000
$ws.Range("I62").Value = 2000
$ws.Range("K62").Value = 2000
$ws.Range("M62").Value = -1376
$ws.Range("H65").Value = 2000
$ws.Range("I65").Value = 2000
$ws.Range("K65").Value = 10000
$ws.Range("M65").Value = -6880
$ws.Range("H115").Value = 25000376
$ws.Range("I115").Value = 33333666
$ws.Range("J115").Value = 500
$ws.Range("K115").Value = 100000998
$ws.Range("L115").Value = 1500
$ws.Range("M115").Value = -99999431
$ws.Range("N115").Value = -4634
$ws.Range("H116").Value = 15194.9
$ws.Range("I116").Value = 36650
$ws.Range("J116").Value = 5999.857
$ws.Range("K116").Value = 36650
$ws.Range("L116").Value = 5999.857
$ws.Range("M116").Value = -33208
$ws.Range("N116").Value = -12883.857
$ws.Range("H127").Value = 1095.7
$ws.Range("I127").Value = 1746.1538
$ws.Range("K127").Value = 5238.4614
$ws.Range("M127").Value = -278.4614000000001
$ws.Range("H135").Value = 487.10526
$ws.Range("J135").Value = 100
$ws.Range("L135").Value = 900
$ws.Range("N135").Value = -5970
$ws.Range("H141").Value = 3104.2354
$ws.Range("I141").Value = 1163.6364
$ws.Range("K141").Value = 3490.9092
$ws.Range("M141").Value = 1689.0908

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1537.579
$ws.Range("I45").Value = 1274.7273
$ws.Range("K45").Value = 1274.7273
$ws.Range("M45").Value = -897.7273
$ws.Range("H74").Value = 1903
$ws.Range("I74").Value = 1938.125
$ws.Range("K74").Value = 1938.125
$ws.Range("M74").Value = -1064.125
$ws.Range("H77").Value = 1903
$ws.Range("I77").Value = 1938.125
$ws.Range("K77").Value = 9690.625
$ws.Range("M77").Value = -5322.625
$ws.Range("H132").Value = 1680.1364
$ws.Range("I132").Value = 1007.3
$ws.Range("J132").Value = 3121.9285
$ws.Range("K132").Value = 3021.9
$ws.Range("L132").Value = 9365.7855
$ws.Range("M132").Value = -491.8999999999996
$ws.Range("N132").Value = -14425.7855
$ws.Range("H138").Value = 25000
$ws.Range("J138").Value = 25000
$ws.Range("L138").Value = 25000
$ws.Range("N138").Value = -35280

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2416.2917
$ws.Range("I105").Value = 2299.6
$ws.Range("J105").Value = 2999.75
$ws.Range("K105").Value = 2299.6
$ws.Range("L105").Value = 2999.75
$ws.Range("M105").Value = -552.5999999999999
$ws.Range("N105").Value = -6493.75
$ws.Range("H107").Value = 2559.2222
$ws.Range("I107").Value = 1861.8572
$ws.Range("K107").Value = 1861.8572
$ws.Range("M107").Value = 58.14280000000008
$ws.Range("H134").Value = 4718.0215
$ws.Range("I134").Value = 5198.5674
$ws.Range("J134").Value = 2940
$ws.Range("K134").Value = 15595.7022
$ws.Range("L134").Value = 8820
$ws.Range("M134").Value = -13060.7022
$ws.Range("N134").Value = -13890

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2190.1428
$ws.Range("I31").Value = 2086.0908
$ws.Range("K31").Value = 2086.0908
$ws.Range("M31").Value = -1791.0908
$ws.Range("H34").Value = 2190.1428
$ws.Range("I34").Value = 2086.0908
$ws.Range("K34").Value = 2086.0908
$ws.Range("M34").Value = -1884.0908
$ws.Range("H45").Value = 7598
$ws.Range("I45").Value = 7598
$ws.Range("K45").Value = 7598
$ws.Range("M45").Value = -7005
$ws.Range("H86").Value = 76925350
$ws.Range("I86").Value = 100002030
$ws.Range("J86").Value = 3098.3333
$ws.Range("K86").Value = 100002030
$ws.Range("L86").Value = 3098.3333
$ws.Range("M86").Value = -100000907
$ws.Range("N86").Value = -5344.3333
$ws.Range("H89").Value = 76925350
$ws.Range("I89").Value = 100002030
$ws.Range("J89").Value = 3098.3333
$ws.Range("K89").Value = 500010150
$ws.Range("L89").Value = 15491.6665
$ws.Range("M89").Value = -500004534
$ws.Range("N89").Value = -26723.6665
$ws.Range("H132").Value = 2249.9697
$ws.Range("J132").Value = 3259.6
$ws.Range("L132").Value = 9778.799999999999
$ws.Range("N132").Value = -14838.8

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 700
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H33").Value = 118.5
$ws.Range("I33").Value = 102.71429
$ws.Range("J33").Value = 229
$ws.Range("K33").Value = 616.28574
$ws.Range("L33").Value = 1374
$ws.Range("M33").Value = -333.28574
$ws.Range("N33").Value = -1940
$ws.Range("H107").Value = 968.6429000000001
$ws.Range("I107").Value = 833.3333
$ws.Range("J107").Value = 1005.5455
$ws.Range("K107").Value = 2499.9999
$ws.Range("L107").Value = 3016.6365
$ws.Range("M107").Value = -579.9998999999998
$ws.Range("N107").Value = -6856.6365
$ws.Range("H136").Value = 125002220
$ws.Range("I136").Value = 125002220
$ws.Range("K136").Value = 375006660
$ws.Range("M136").Value = -375001560

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1325
$ws.Range("J113").Value = 1325
$ws.Range("L113").Value = 1325
$ws.Range("N113").Value = -5665
$ws.Range("H122").Value = 2449.5
$ws.Range("J122").Value = 2500
$ws.Range("L122").Value = 7500
$ws.Range("N122").Value = -12400
$ws.Range("H136").Value = 10602.454
$ws.Range("J136").Value = 10602.454
$ws.Range("L136").Value = 31807.362
$ws.Range("N136").Value = -36907.362
$ws.Range("H140").Value = 38539.855
$ws.Range("J140").Value = 38539.855
$ws.Range("L140").Value = 38539.855
$ws.Range("N140").Value = -48899.855

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4280
$ws.Range("I22").Value = 6200
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 6200
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = -5905
$ws.Range("N22").Value = -3590
$ws.Range("H27").Value = 4280
$ws.Range("I27").Value = 6200
$ws.Range("J27").Value = 3000
$ws.Range("K27").Value = 6200
$ws.Range("L27").Value = 3000
$ws.Range("M27").Value = -6093
$ws.Range("N27").Value = -3214
$ws.Range("H42").Value = 29054.938
$ws.Range("J42").Value = 29054.938
$ws.Range("L42").Value = 29054.938
$ws.Range("N42").Value = -30180.938
$ws.Range("H49").Value = 29054.938
$ws.Range("J49").Value = 29054.938
$ws.Range("L49").Value = 29054.938
$ws.Range("N49").Value = -29348.938

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H140").Value = 55117.25
$ws.Range("J140").Value = 55117.25
$ws.Range("L140").Value = 55117.25
$ws.Range("N140").Value = -65477.25
$ws.Range("H141").Value = 80525.234
$ws.Range("J141").Value = 80525.234
$ws.Range("L141").Value = 80525.234
$ws.Range("N141").Value = -90885.234
